# Daily "Updated cryptos list" refresh (GitHub Actions job).
# Price (column D) and Volume(1h) (column E) cells are stored as literal
# text in the sheet (periods used as thousands separators, e.g. "63.298.05"),
# so every write below forces text storage - even for values that look like
# plain numbers ("0.999", "588.26", ...) - by flipping the cell to the "@"
# (Text) number format before the write and restoring the default "Normal"
# cell style right after, so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value2 = '63.298.05'
$ws.Range("E2").Value2 = '  +0.19%  '

$ws.Range("D3").Value2 = '2.595.12'
$ws.Range("E3").Value2 = '  +1.60%  '

Set-TextValue "D4" '0.999'
$ws.Range("E4").Value2 = '  -0.05%  '

Set-TextValue "D5" '588.26'
$ws.Range("E5").Value2 = '  +3.46%  '

Set-TextValue "D6" '149.06'
$ws.Range("E6").Value2 = '  +1.54%  '

Set-TextValue "D7" '0.999'
$ws.Range("E7").Value2 = '  -0.03%  '

$ws.Range("E8").Value2 = '  +1.92%  '

$ws.Range("E9").Value2 = '  +4.05%  '

$ws.Range("E10").Value2 = '  +1.65%  '

$ws.Range("E11").Value2 = '  +0.00%  '

$ws.Range("E12").Value2 = '  +1.48%  '

Set-TextValue "D13" '27.72'
$ws.Range("E13").Value2 = '  +0.68%  '

$ws.Range("D14").Value2 = '3.054.92'
$ws.Range("E14").Value2 = '  +1.65%  '

$ws.Range("D15").Value2 = '63.278.01'
$ws.Range("E15").Value2 = '  +0.28%  '

$ws.Range("D17").Value2 = '2.572.80'
$ws.Range("E17").Value2 = '  +0.76%  '

$ws.Range("E18").Value2 = '  +0.12%  '

Set-TextValue "D19" '345.47'
$ws.Range("E19").Value2 = '  +2.95%  '

$ws.Range("E20").Value2 = '  +2.84%  '

$ws.Range("E21").Value2 = '  +1.53%  '

$ws.Range("E22").Value2 = '  +0.00%  '

$ws.Range("E23").Value2 = '  -3.40%  '

Set-TextValue "D24" '66.86'
$ws.Range("E24").Value2 = '  +2.46%  '

$ws.Range("E25").Value2 = '  +0.48%  '

$ws.Range("D26").Value2 = '2.670.14'
$ws.Range("E26").Value2 = '  -0.21%  '

$ws.Range("E27").Value2 = '  -0.57%  '

Set-TextValue "D28" '8.26'
$ws.Range("E28").Value2 = '  +12.40%  '

Set-TextValue "D29" '8.52'
$ws.Range("E29").Value2 = '  +0.52%  '

$ws.Range("B30").Value2 = 'SuiNetwork'
$ws.Range("C30").Value2 = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D30" '1.49'
$ws.Range("E30").Value2 = '  +0.78%  '

$ws.Range("B31").Value2 = 'Binance-PegBSC-USD'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D31" '0.999'
$ws.Range("E31").Value2 = '  -0.10%  '

Set-TextValue "D32" '2.01'
$ws.Range("E32").Value2 = '  +8.21%  '

$ws.Range("D33").Value2 = '0.0₃0833'
$ws.Range("E33").Value2 = '  +1.52%  '

Set-TextValue "D34" '467.88'
$ws.Range("E34").Value2 = '  +14.60%  '

$ws.Range("E35").Value2 = '  +4.33%  '

Set-TextValue "D36" '177.22'
$ws.Range("E36").Value2 = '  +0.72%  '

$ws.Range("E37").Value2 = '  +2.35%  '

$ws.Range("E38").Value2 = '  +1.68%  '

$ws.Range("E39").Value2 = '  +6.40%  '

$ws.Range("E40").Value2 = '  +0.04%  '

$ws.Range("E41").Value2 = '  +0.49%  '

$ws.Range("E42").Value2 = '  +0.01%  '

Set-TextValue "D43" '152.49'
$ws.Range("E43").Value2 = '  -0.28%  '

$ws.Range("E44").Value2 = '  +1.98%  '

$ws.Range("E45").Value2 = '  +0.43%  '

$ws.Range("E46").Value2 = '  +5.28%  '

Set-TextValue "D47" '0.617'
$ws.Range("E47").Value2 = '  +1.78%  '

$ws.Range("E48").Value2 = '  +1.66%  '

$ws.Range("E49").Value2 = '  +1.27%  '

$ws.Range("E50").Value2 = '  -0.58%  '

$ws.Range("E51").Value2 = '  +0.85%  '
